$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B80").Value = 5498503
$ws.Range("F80").Value = "FC Koper"
$ws.Range("G80").Value = "NS Mura"
$ws.Range("H80").Value = 1
$ws.Range("I80").Value = 2
$ws.Range("J80").Value = "A"
$ws.Range("K80").Value = 2.05
$ws.Range("L80").Value = 3.3
$ws.Range("M80").Value = 3.25
$ws.Range("N80").Value = 2
$ws.Range("O80").Value = 3.4
$ws.Range("P80").Value = 3.25
$ws.Range("Q80").Value = -0.5
$ws.Range("R80").Value = 2
$ws.Range("S80").Value = 1.8
$ws.Range("T80").Value = 2.5
$ws.Range("U80").Value = 1.825
$ws.Range("V80").Value = 1.975
$ws.Range("W80").Value = -1
$ws.Range("X80").Value = -1
$ws.Range("Y80").Value = 2.25
$ws.Range("Z80").Value = -1
$ws.Range("AA80").Value = 0.8
$ws.Range("AB80").Value = 0.825
$ws.Range("AC80").Value = -1

$ws.Range("B81").Value = 5495053
$ws.Range("F81").Value = "NK Radomlje"
$ws.Range("G81").Value = "NK Domzale"
$ws.Range("H81").Value = 1
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = "H"
$ws.Range("K81").Value = 2.55
$ws.Range("L81").Value = 3.1
$ws.Range("M81").Value = 2.55
$ws.Range("N81").Value = 3.75
$ws.Range("O81").Value = 3.4
$ws.Range("P81").Value = 1.833
$ws.Range("Q81").Value = 0.5
$ws.Range("R81").Value = 1.925
$ws.Range("S81").Value = 1.875
$ws.Range("T81").Value = 2.5
$ws.Range("U81").Value = 1.975
$ws.Range("V81").Value = 1.825
$ws.Range("W81").Value = 2.75
$ws.Range("X81").Value = -1
$ws.Range("Y81").Value = -1
$ws.Range("Z81").Value = 0.925
$ws.Range("AA81").Value = -1
$ws.Range("AB81").Value = -1
$ws.Range("AC81").Value = 0.825

$ws.Range("B82").Value = 6816473
$ws.Range("F82").Value = "NK Bravo"
$ws.Range("G82").Value = "NK Rogaska"
$ws.Range("H82").Value = 2
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = "H"
$ws.Range("K82").Value = 1.8
$ws.Range("L82").Value = 3.5
$ws.Range("M82").Value = 4
$ws.Range("N82").Value = 2.05
$ws.Range("O82").Value = 3
$ws.Range("P82").Value = 3.75
$ws.Range("Q82").Value = -0.25
$ws.Range("R82").Value = 1.75
$ws.Range("S82").Value = 2.05
$ws.Range("T82").Value = 2.25
$ws.Range("U82").Value = 1.95
$ws.Range("V82").Value = 1.85
$ws.Range("W82").Value = 1.05
$ws.Range("X82").Value = -1
$ws.Range("Y82").Value = -1
$ws.Range("Z82").Value = 0.75
$ws.Range("AA82").Value = -1
$ws.Range("AB82").Value = -0.5
$ws.Range("AC82").Value = 0.425

$ws.Range("B83").Value = 6814327
$ws.Range("F83").Value = "NS Mura"
$ws.Range("G83").Value = "NK Domzale"
$ws.Range("H83").Value = 2
$ws.Range("I83").Value = 3
$ws.Range("J83").Value = "A"
$ws.Range("K83").Value = 2
$ws.Range("L83").Value = 3.3
$ws.Range("M83").Value = 3.4
$ws.Range("N83").Value = 1.909
$ws.Range("O83").Value = 3.4
$ws.Range("P83").Value = 3.75
$ws.Range("Q83").Value = -0.5
$ws.Range("R83").Value = 1.95
$ws.Range("S83").Value = 1.85
$ws.Range("T83").Value = 2.5
$ws.Range("U83").Value = 1.9
$ws.Range("V83").Value = 1.9
$ws.Range("W83").Value = -1
$ws.Range("X83").Value = -1
$ws.Range("Y83").Value = 2.75
$ws.Range("Z83").Value = -1
$ws.Range("AA83").Value = 0.8500000000000001
$ws.Range("AB83").Value = 0.8999999999999999
$ws.Range("AC83").Value = -1

$ws.Range("B181").Value = 7680773
$ws.Range("E181").Value = 45340.375
$ws.Range("F181").Value = "NK Bravo"
$ws.Range("G181").Value = "NK Celje"
$ws.Range("K181").Value = 4.5
$ws.Range("L181").Value = 3.4
$ws.Range("M181").Value = 1.75
$ws.Range("N181").Value = 5
$ws.Range("O181").Value = 3.6
$ws.Range("P181").Value = 1.7
$ws.Range("Q181").Value = 0.75
$ws.Range("R181").Value = 1.85
$ws.Range("S181").Value = 1.95
$ws.Range("T181").Value = 2.5
$ws.Range("U181").Value = 1.975
$ws.Range("V181").Value = 1.825

$ws.Range("B182").Value = 7680776
$ws.Range("E182").Value = 45340.45833333334
$ws.Range("F182").Value = "NS Mura"
$ws.Range("G182").Value = "NK Maribor"
$ws.Range("K182").Value = 3.2
$ws.Range("L182").Value = 3.2
$ws.Range("M182").Value = 2.15
$ws.Range("N182").Value = 3.8
$ws.Range("O182").Value = 3.3
$ws.Range("P182").Value = 1.833
$ws.Range("Q182").Value = 0.5
$ws.Range("R182").Value = 1.975
$ws.Range("S182").Value = 1.825
$ws.Range("T182").Value = 2.75
$ws.Range("U182").Value = 1.975
$ws.Range("V182").Value = 1.825

# Remove the now-obsolete last row (old row 183), shifting dimension to AC182
$ws.Rows(183).Delete() | Out-Null
